# Update countries & provincias Spain
#
# The "Pais" sheet lists one country per row (col A) with daily COVID-19
# stats in columns B-H (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes). This update:
#   - refreshes the case/death counters for several existing countries
#   - inserts several newly-tracked countries (Luxemburgo, Jordania,
#     Mayotte, Aruba, San Martin (Parte Francesa), Haiti, Birmania,
#     Surinam, Benin, Santa Sede, Nepal, Sudan, Cabo Verde) ahead of
#     where they used to sit in the list, which shifts the countries
#     that used to occupy those rows down by one position
#
# Because the row position (not the country name) is what drives each
# cell address, the edit is expressed as direct per-cell writes: for
# every affected row we set the country name (col A, only where it
# changed) together with the final B:H values for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B4').Value = 113677
$ws.Range('C4').Value = 9551
$ws.Range('E4').Value = 108555
$ws.Range('G4').Value = 207
$ws.Range('H4').Value = 1903

$ws.Range('B5').Value = 92472
$ws.Range('C5').Value = 5974
$ws.Range('D5').Value = 12384
$ws.Range('E5').Value = 70065
$ws.Range('F5').Value = 3856
$ws.Range('G5').Value = 889
$ws.Range('H5').Value = 10023

$ws.Range('B8').Value = 56202
$ws.Range('C8').Value = 5331
$ws.Range('E8').Value = 49141
$ws.Range('G8').Value = 52
$ws.Range('H8').Value = 403

$ws.Range('B16').Value = 8122
$ws.Range('C16').Value = 425
$ws.Range('E16').Value = 7829

$ws.Range('B17').Value = 7402
$ws.Range('C17').Value = 1704
$ws.Range('E17').Value = 7252
$ws.Range('F17').Value = 309
$ws.Range('G17').Value = 16
$ws.Range('H17').Value = 108

$ws.Range('B19').Value = 4933
$ws.Range('C19').Value = 176
$ws.Range('E19').Value = 4524

$ws.Range('B24').Value = 3447
$ws.Range('C24').Value = 378
$ws.Range('E24').Value = 3326
$ws.Range('F24').Value = 239

$ws.Range('A30').Value = 'Luxemburgo'
$ws.Range('B30').Value = 1831
$ws.Range('C30').Value = 226
$ws.Range('D30').Value = 40
$ws.Range('E30').Value = 1773
$ws.Range('F30').Value = 25
$ws.Range('G30').Value = 3
$ws.Range('H30').Value = 18

$ws.Range('A31').Value = 'Ecuador'
$ws.Range('B31').Value = 1823
$ws.Range('C31').Value = 196
$ws.Range('D31').Value = 3
$ws.Range('E31').Value = 1772
$ws.Range('F31').Value = 58
$ws.Range('G31').Value = 7
$ws.Range('H31').Value = 48

$ws.Range('D44').Value = 114
$ws.Range('E44').Value = 847

$ws.Range('B57').Value = 590
$ws.Range('C57').Value = 28
$ws.Range('D57').Value = 45
$ws.Range('E57').Value = 544
$ws.Range('G57').Value = 1
$ws.Range('H57').Value = 1

$ws.Range('B79').Value = 258
$ws.Range('C79').Value = 21
$ws.Range('G79').Value = 1
$ws.Range('H79').Value = 5

$ws.Range('A80').Value = 'Jordania'
$ws.Range('B80').Value = 246
$ws.Range('C80').Value = 11
$ws.Range('D80').Value = 18
$ws.Range('E80').Value = 227
$ws.Range('F80').Value = 3
$ws.Range('G80').Value = 0
$ws.Range('H80').Value = 1

$ws.Range('A81').Value = 'Republica de Macedonia'
$ws.Range('B81').Value = 241
$ws.Range('C81').Value = 22
$ws.Range('D81').Value = 3
$ws.Range('E81').Value = 234
$ws.Range('F81').Value = 1
$ws.Range('G81').Value = 1
$ws.Range('H81').Value = 4

$ws.Range('A118').Value = 'Mayotte'
$ws.Range('B118').Value = 63
$ws.Range('C118').Value = 13
$ws.Range('E118').Value = 63

$ws.Range('A119').Value = 'Kirguistan'
$ws.Range('C119').Value = 0
$ws.Range('D119').Value = 0
$ws.Range('E119').Value = 58
$ws.Range('G119').Value = 0
$ws.Range('H119').Value = 0

$ws.Range('A120').Value = 'Consejo Danes para los Refugiados'
$ws.Range('B120').Value = 58
$ws.Range('C120').Value = 7
$ws.Range('D120').Value = 2
$ws.Range('E120').Value = 50
$ws.Range('G120').Value = 3
$ws.Range('H120').Value = 6

$ws.Range('A121').Value = 'Liechtenstein'
$ws.Range('C121').Value = 0
$ws.Range('D121').Value = 0
$ws.Range('E121').Value = 56
$ws.Range('F121').Value = 0
$ws.Range('H121').Value = 0

$ws.Range('A122').Value = 'Paraguay'
$ws.Range('C122').Value = 4
$ws.Range('D122').Value = 1
$ws.Range('E122').Value = 52
$ws.Range('F122').Value = 1
$ws.Range('H122').Value = 3

$ws.Range('A123').Value = 'Gibraltar'
$ws.Range('B123').Value = 56
$ws.Range('C123').Value = 1
$ws.Range('D123').Value = 14
$ws.Range('E123').Value = 42

$ws.Range('A124').Value = 'Ruanda'
$ws.Range('B124').Value = 54
$ws.Range('E124').Value = 54

$ws.Range('A127').Value = 'Aruba'
$ws.Range('B127').Value = 40
$ws.Range('C127').Value = 7
$ws.Range('E127').Value = 39
$ws.Range('H127').Value = 0

$ws.Range('A128').Value = 'Puerto Rico'
$ws.Range('B128').Value = 39
$ws.Range('C128').Value = 0
$ws.Range('F128').Value = 0
$ws.Range('H128').Value = 2

$ws.Range('A129').Value = 'Kenia'
$ws.Range('B129').Value = 38
$ws.Range('C129').Value = 7
$ws.Range('D129').Value = 1
$ws.Range('E129').Value = 36
$ws.Range('F129').Value = 2
$ws.Range('H129').Value = 1

$ws.Range('A130').Value = 'Macao'
$ws.Range('B130').Value = 34
$ws.Range('D130').Value = 10
$ws.Range('E130').Value = 24

$ws.Range('A153').Value = 'San Martin (Parte Francesa)'

$ws.Range('A154').Value = 'Dominica'

$ws.Range('A155').Value = 'Niger'
$ws.Range('D155').Value = 0
$ws.Range('H155').Value = 1

$ws.Range('A156').Value = 'Bahamas'
$ws.Range('D156').Value = 1
$ws.Range('H156').Value = 0

$ws.Range('A159').Value = 'Haiti'
$ws.Range('C159').Value = 0

$ws.Range('A160').Value = 'Birmania'

$ws.Range('A161').Value = 'Surinam'

$ws.Range('A162').Value = 'Mozambique'
$ws.Range('C162').Value = 1

$ws.Range('A172').Value = 'Benin'

$ws.Range('A173').Value = 'Laos'

$ws.Range('A174').Value = 'Santa Sede'
$ws.Range('C174').Value = 2

$ws.Range('A175').Value = 'Eritrea'
$ws.Range('B175').Value = 6
$ws.Range('E175').Value = 6

$ws.Range('A176').Value = 'San Bartolome'

$ws.Range('A177').Value = 'Fiyi'

$ws.Range('A178').Value = 'Montserrat'

$ws.Range('A179').Value = 'Siria'
$ws.Range('C179').Value = 0

$ws.Range('A180').Value = 'Mauritania'
$ws.Range('C180').Value = 2
$ws.Range('E180').Value = 5
$ws.Range('H180').Value = 0

$ws.Range('A184').Value = 'Guyana'
$ws.Range('B184').Value = 5
$ws.Range('H184').Value = 1
